$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "41.080.81"
Set-TextValue "E2" "  -1.36%  "

Set-TextValue "D3" "2.424.84"
Set-TextValue "E3" "  -1.98%  "

Set-TextValue "E4" "  +0.21%  "

Set-TextValue "D5" "316.61"
Set-TextValue "E5" "  -0.39%  "

Set-TextValue "D6" "89.02"
Set-TextValue "E6" "  -3.99%  "

Set-TextValue "D7" "0.540"
Set-TextValue "E7" "  -2.44%  "

Set-TextValue "E8" "  +0.13%  "

Set-TextValue "E9" "  -4.52%  "

Set-TextValue "B10" "Avalanche"
Set-TextValue "C10" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D10" "31.95"
Set-TextValue "E10" "  -3.29%  "

Set-TextValue "B11" "Dogecoin"
Set-TextValue "C11" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D11" "0.0831"
Set-TextValue "E11" "  -4.23%  "

Set-TextValue "E12" "  -2.62%  "

Set-TextValue "D13" "2.801.07"
Set-TextValue "E13" "  -1.84%  "

Set-TextValue "D14" "6.69"
Set-TextValue "E14" "  -3.18%  "

Set-TextValue "D15" "15.54"
Set-TextValue "E15" "  -0.38%  "

Set-TextValue "D16" "2.416.85"
Set-TextValue "E16" "  -2.30%  "

Set-TextValue "D17" "0.770"
Set-TextValue "E17" "  -2.47%  "

Set-TextValue "D18" "41.005.31"
Set-TextValue "E18" "  -1.42%  "

Set-TextValue "D19" "0.0₃0921"
Set-TextValue "E19" "  -3.55%  "

Set-TextValue "D20" "6.24"
Set-TextValue "E20" "  -3.97%  "

Set-TextValue "D21" "71.88"
Set-TextValue "E21" "  +0.75%  "

Set-TextValue "D22" "11.03"
Set-TextValue "E22" "  -3.07%  "

Set-TextValue "D23" "234.80"
Set-TextValue "E23" "  -2.40%  "

Set-TextValue "D24" "2.69"
Set-TextValue "E24" "  -2.01%  "

Set-TextValue "E25" "  +0.08%  "

Set-TextValue "D26" "1.87"
Set-TextValue "E26" "  -3.14%  "

Set-TextValue "D27" "24.08"
Set-TextValue "E27" "  -2.61%  "

Set-TextValue "E28" "  -3.48%  "

Set-TextValue "E29" "  -4.02%  "

Set-TextValue "D30" "34.70"
Set-TextValue "E30" "  -5.06%  "

Set-TextValue "D31" "156.24"
Set-TextValue "E31" "  -1.73%  "

Set-TextValue "E32" "  +0.02%  "

Set-TextValue "D33" "5.25"
Set-TextValue "E33" "  -5.42%  "

Set-TextValue "B34" "WEMIXToken"
Set-TextValue "C34" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D34" "2.51"
Set-TextValue "E34" "  -2.72%  "

Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.0744"
Set-TextValue "E35" "  -3.74%  "

Set-TextValue "E36" "  +0.78%  "

Set-TextValue "D37" "16.64"
Set-TextValue "E37" "  -4.58%  "

Set-TextValue "E38" "  -0.92%  "

Set-TextValue "D39" "1.77"
Set-TextValue "E39" "  -3.84%  "

Set-TextValue "E40" "  -2.96%  "

Set-TextValue "E41" "  -3.32%  "

Set-TextValue "E42" "  -7.58%  "

Set-TextValue "D43" "1.985.75"
Set-TextValue "E43" "  -0.02%  "

Set-TextValue "D44" "18.67"
Set-TextValue "E44" "  -2.90%  "

Set-TextValue "E45" "  -3.73%  "

Set-TextValue "D46" "2.87"
Set-TextValue "E46" "  -4.92%  "

Set-TextValue "D47" "9.45"
Set-TextValue "E47" "  +2.16%  "

Set-TextValue "D48" "2.659.95"
Set-TextValue "E48" "  -1.87%  "

Set-TextValue "D49" "94.86"
Set-TextValue "E49" "  -2.81%  "

Set-TextValue "D50" "72.97"
Set-TextValue "E50" "  -0.82%  "

Set-TextValue "D51" "51.77"
Set-TextValue "E51" "  -1.71%  "

